$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '25.890.20'
$ws.Range('E2').Value = '  +0.21%  '

# Row 3
Set-TextValue $ws.Range('D3') '1.646.45'
$ws.Range('E3').Value = '  +0.76%  '

# Row 4
Set-TextValue $ws.Range('D4') '1.008'
$ws.Range('E4').Value = '  +0.62%  '

# Row 5
Set-TextValue $ws.Range('D5') '215.01'
$ws.Range('E5').Value = '  +0.06%  '

# Row 6
Set-TextValue $ws.Range('D6') '0.5064'

# Row 7
$ws.Range('E7').Value = '  +0.42%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.2572'
$ws.Range('E8').Value = '  -0.09%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.06407'
$ws.Range('E9').Value = '  -0.02%  '

# Row 10
Set-TextValue $ws.Range('D10') '19.73'
$ws.Range('E10').Value = '  +0.49%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.07767'
$ws.Range('E11').Value = '  +1.33%  '

# Row 12
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D12') '4.306'
$ws.Range('E12').Value = '  +1.59%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D13') '1.662.11'
$ws.Range('E13').Value = '  +1.69%  '

# Row 14
$ws.Range('E14').Value = '  -0.12%  '

# Row 15
Set-TextValue $ws.Range('D15') '0.0₅7889'
$ws.Range('E15').Value = '  -0.48%  '

# Row 16
Set-TextValue $ws.Range('D16') '65.14'
$ws.Range('E16').Value = '  +2.57%  '

# Row 17
Set-TextValue $ws.Range('D17') '25.970.53'
$ws.Range('E17').Value = '  +0.55%  '

# Row 18
Set-TextValue $ws.Range('D18') '1.007'
$ws.Range('E18').Value = '  +0.46%  '

# Row 19
Set-TextValue $ws.Range('D19') '197.30'
$ws.Range('E19').Value = '  -2.79%  '

# Row 20
$ws.Range('E20').Value = '  +2.25%  '

# Row 21
$ws.Range('E21').Value = '  +0.68%  '

# Row 22
Set-TextValue $ws.Range('D22') '6.059'
$ws.Range('E22').Value = '  +1.61%  '

# Row 23
$ws.Range('E23').Value = '  +0.54%  '

# Row 24
Set-TextValue $ws.Range('D24') '1.864'
$ws.Range('E24').Value = '  -3.64%  '

# Row 25
Set-TextValue $ws.Range('D25') '141.11'
$ws.Range('E25').Value = '  +0.15%  '

# Row 26
$ws.Range('E26').Value = '  +0.36%  '

# Row 27
Set-TextValue $ws.Range('D27') '6.895'
$ws.Range('E27').Value = '  +3.04%  '

# Row 28
Set-TextValue $ws.Range('D28') '15.71'
$ws.Range('E28').Value = '  +0.02%  '

# Row 29
Set-TextValue $ws.Range('D29') '0.05075'
$ws.Range('E29').Value = '  +1.97%  '

# Row 30
Set-TextValue $ws.Range('D30') '1.240'
$ws.Range('E30').Value = '  +0.03%  '

# Row 31
Set-TextValue $ws.Range('D31') '3.267'
$ws.Range('E31').Value = '  -0.23%  '

# Row 32
Set-TextValue $ws.Range('D32') '3.202'
$ws.Range('E32').Value = '  +0.60%  '

# Row 33
Set-TextValue $ws.Range('D33') '1.543'
$ws.Range('E33').Value = '  +0.55%  '

# Row 34
Set-TextValue $ws.Range('D34') '2.370'
$ws.Range('E34').Value = '  +0.79%  '

# Row 35
Set-TextValue $ws.Range('D35') '0.8943'
$ws.Range('E35').Value = '  +0.23%  '

# Row 36
Set-TextValue $ws.Range('D36') '2.598'
$ws.Range('E36').Value = '  -0.66%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.5545'

# Row 38
Set-TextValue $ws.Range('D38') '1.131.62'
$ws.Range('E38').Value = '  -3.57%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.01563'
$ws.Range('E39').Value = '  +0.42%  '

# Row 40
Set-TextValue $ws.Range('D40') '1.007'
$ws.Range('E40').Value = '  +0.56%  '

# Row 41
Set-TextValue $ws.Range('D41') '5.678'
$ws.Range('E41').Value = '  +0.53%  '

# Row 42
Set-TextValue $ws.Range('D42') '0.8154'
$ws.Range('E42').Value = '  +1.48%  '

# Row 43
Set-TextValue $ws.Range('D43') '99.71'
$ws.Range('E43').Value = '  +0.35%  '

# Row 44
$ws.Range('E44').Value = '  +6.77%  '

# Row 45
Set-TextValue $ws.Range('D45') '1.780.50'
$ws.Range('E45').Value = '  +0.58%  '

# Row 46
Set-TextValue $ws.Range('D46') '0.4542'
$ws.Range('E46').Value = '  +0.70%  '

# Row 47
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D47') '55.28'
$ws.Range('E47').Value = '  +0.87%  '

# Row 48
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range('D48') '1.006'
$ws.Range('E48').Value = '  +0.05%  '

# Row 49
Set-TextValue $ws.Range('D49') '0.05090'
$ws.Range('E49').Value = '  +1.20%  '

# Row 50
$ws.Range('E50').Value = '  +0.56%  '

# Row 51
Set-TextValue $ws.Range('D51') '0.09563'
$ws.Range('E51').Value = '  +3.08%  '
